$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 171
$ws1.Range("F4").Value = 172
$ws1.Range("F5").Value = 4897
$ws1.Range("F9").Value = 536
$ws1.Range("F10").Value = 493
$ws1.Range("F13").Value = 1364
$ws1.Range("F14").Value = 3359
$ws1.Range("F15").Value = 397
$ws1.Range("F16").Value = 123
$ws1.Range("F17").Value = 108
$ws1.Range("F18").Value = 72
$ws1.Range("F19").Value = 2535
$ws1.Range("F20").Value = 124
$ws1.Range("F21").Value = 81
$ws1.Range("F24").Value = 38
$ws1.Range("F25").Value = 122
$ws1.Range("F26").Value = 54

# Sheet "全部类型" (fourth sheet) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 171
$ws4.Range("F4").Value = 172
$ws4.Range("F6").Value = 4897
$ws4.Range("F10").Value = 536
$ws4.Range("F11").Value = 493
$ws4.Range("F14").Value = 1364
$ws4.Range("F15").Value = 3359
$ws4.Range("F16").Value = 397
$ws4.Range("F17").Value = 123
$ws4.Range("F18").Value = 108
$ws4.Range("F19").Value = 72
$ws4.Range("F20").Value = 2535
$ws4.Range("F21").Value = 124
$ws4.Range("F22").Value = 81
$ws4.Range("F25").Value = 38
$ws4.Range("F26").Value = 122
$ws4.Range("F27").Value = 54
